$wb = $excel.ActiveWorkbook
$src = $wb.Worksheets.Item("6")
$count = $wb.Worksheets.Count
$last = $wb.Worksheets.Item($count)
$src.Copy($null, $last)
$newSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$newSheet.Name = "11"
foreach ($s in $wb.Worksheets) {
    Write-Host $s.Name
}
